$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Topics")

# Update the "No. of Hrs" value for row 2 (first topic, "Introduction")
# from 30mins to 20mins
$ws.Range("D2").Value = "20mins"

# Add the new sub-topic "Verify SSH client" under "Workstation Setup" (row 9)
$ws.Range("C9").Value = "Verify SSH client"

# Move/update the active selection to D2 (matches the saved selection state)
$ws.Range("D2").Select()
